# Trim trailing whitespace from the company-name cells in column C
# across all 4 "Turno" sheets, as described in the diff:
#   "CORPORACIÓN ENRIQUE & JANIS " -> "CORPORACIÓN ENRIQUE & JANIS"
#   "MULTITRANS RR "               -> "MULTITRANS RR"
# These occur in rows 7,8,9,10,36,37,38,39 (CORPORACIÓN ...) and
# rows 93,94,95,96 (MULTITRANS RR) on every sheet.

$wb = $excel.ActiveWorkbook

$rowsCorp = @(7, 8, 9, 10, 36, 37, 38, 39)
$rowsMulti = @(93, 94, 95, 96)

$allRows = $rowsCorp + $rowsMulti

foreach ($ws in $wb.Worksheets) {
    foreach ($r in $allRows) {
        $cell = $ws.Cells.Item($r, 3)  # Column C
        $current = $cell.Value()
        if ($current -ne $null) {
            $trimmed = ([string]$current).TrimEnd()
            if ($trimmed -ne $current) {
                $cell.Value = $trimmed
            }
        }
    }
}
